# Apply updated CircadiPy cosinor analysis results (sawtooth_0.5 simulation)
# to rows 2 and 3 of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 -----------------------------------------------------------
$ws.Range("E2").Value = 25.99000000000062

$ws.Range("G2").Value = 0.001300756447688678
$ws.Range("H2").Value = 0.004413277499858736

$ws.Range("K2").Value = 4.47052274004391
$ws.Range("L2").Value = "[1.4611416100562193, 7.479903870031602]"

$ws.Range("M2").Value = 0.003713199346918161
$ws.Range("N2").Value = 0.003713199346918161

$ws.Range("O2").Value = -1.094368612061309
$ws.Range("P2").Value = "[-1.8113687372049263, -0.37736848691769254]"

$ws.Range("Q2").Value = 0.002880827789364382
$ws.Range("R2").Value = 0.002880827789364382

$ws.Range("S2").Value = 13.68725948248285
$ws.Range("T2").Value = "[11.994772481895358, 15.379746483070344]"

$ws.Range("W2").Value = 4.526786786786896
$ws.Range("X2").Value = 1.560960960960998
$ws.Range("Y2").Value = 7.492612612612795

# --- Row 3 -----------------------------------------------------------
$ws.Range("E3").Value = 23.9000000000003

$ws.Range("G3").Value = 0.0002286941940167209
$ws.Range("H3").Value = 0.002572995977249508
$ws.Range("I3").Value = 0.04174847316225927

$ws.Range("K3").Value = 5.419939239656232
$ws.Range("L3").Value = "[2.1098310164758125, 8.730047462836652]"

$ws.Range("M3").Value = 0.001416241640530336
$ws.Range("N3").Value = 0.002832483281060671

$ws.Range("O3").Value = 1.943447707626118
$ws.Range("P3").Value = "[1.2893423303021159, 2.5975530849501194]"

$ws.Range("Q3").Value = 0.0000000135992404004525708
$ws.Range("R3").Value = 0.0000000271984808009051416

$ws.Range("S3").Value = 13.52266845877242
$ws.Range("T3").Value = "[11.708976403655761, 15.336360513889087]"

$ws.Range("W3").Value = 16.50750750750771
$ws.Range("X3").Value = 14.01941941941959
$ws.Range("Y3").Value = 18.99559559559584
